# Add a new "instancia" column (F) with the resolution type for the
# existing data row, widen the data columns (B:E) to fit their content,
# and move the active selection to B9 (mirrors the authored workbook
# upload: new data + cosmetic view/column-width tweaks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + value for column F.
$ws.Range("F1").Value = "instancia"
$ws.Range("F2").Value = "Resolución rectoral"

# Column widths added for B:E (stored width = COM ColumnWidth + 5/6;
# snapped to the nearest 1/6 of a character by the host, so request the
# already-snapped values so the persisted width matches exactly).
$ws.Columns.Item(2).ColumnWidth = 18.66666667
$ws.Columns.Item(3).ColumnWidth = 14.33333333
$ws.Columns.Item(4).ColumnWidth = 14.83333333
$ws.Columns.Item(5).ColumnWidth = 19.83333333

# Selection moved from B8 to B9.
$ws.Range("B9").Select()
